# Update cfb_weather.xlsx with Timestamp 2024-11-11T10:01:41.103846
# Refreshes the scraped weather/odds snapshot (FD line movement, wind, etc.)
# and bumps every row's Timestamp to the new scrape time.

$wb = $excel.ActiveWorkbook
$wsFBS = $wb.Worksheets.Item("FBS")
$wsOther = $wb.Worksheets.Item("Other")

# --- Weather / odds data corrections (sheet "FBS") ---
$wsFBS.Range("Y2").Value = 49.5
$wsFBS.Range("Z2").Value = -106
$wsFBS.Range("AE2").Value = -0.0198019801980198
$wsFBS.Range("AB3").Value = -15.5
$wsFBS.Range("AF3").Value = -0.5
$wsFBS.Range("Z6").Value = -115
$wsFBS.Range("AB6").Value = 4
$wsFBS.Range("AF6").Value = 0
$wsFBS.Range("Y8").Value = 53.5
$wsFBS.Range("Z8").Value = -110
$wsFBS.Range("AE8").Value = 0.01904761904761905
$wsFBS.Range("AB9").Value = -15.5
$wsFBS.Range("AF9").Value = 1
$wsFBS.Range("Y10").Value = 58.5
$wsFBS.Range("AB10").Value = -9
$wsFBS.Range("AE10").Value = 0.01739130434782609
$wsFBS.Range("AF10").Value = 1
$wsFBS.Range("Y11").Value = 53.5
$wsFBS.Range("Z11").Value = -108
$wsFBS.Range("AE11").Value = 0.01904761904761905
$wsFBS.Range("Z14").Value = -114
$wsFBS.Range("Y15").Value = 48.5
$wsFBS.Range("AE15").Value = 0.02105263157894737
$wsFBS.Range("Y16").Value = 47.5
$wsFBS.Range("Z16").Value = -112
$wsFBS.Range("AE16").Value = 0.04395604395604396
$wsFBS.Range("AB19").Value = -13.5
$wsFBS.Range("AF19").Value = 0.5
$wsFBS.Range("Y20").Value = 54.5
$wsFBS.Range("AB20").Value = 10
$wsFBS.Range("AE20").Value = -0.01801801801801802
$wsFBS.Range("AF20").Value = -1
$wsFBS.Range("Z22").Value = -110
$wsFBS.Range("Z24").Value = -106
$wsFBS.Range("Z25").Value = -105
$wsFBS.Range("Y27").Value = 42.5
$wsFBS.Range("AB27").Value = 14
$wsFBS.Range("AE27").Value = 0.02409638554216868
$wsFBS.Range("AF27").Value = -1
$wsFBS.Range("Y29").Value = 59.5
$wsFBS.Range("AE29").Value = 0.0170940170940171
$wsFBS.Range("Z31").Value = -105
$wsFBS.Range("AB31").Value = 29.5
$wsFBS.Range("AF31").Value = -2
$wsFBS.Range("AB33").Value = -16.5
$wsFBS.Range("AF33").Value = 1.5
$wsFBS.Range("AB35").Value = 5
$wsFBS.Range("AF35").Value = -0.5
$wsFBS.Range("Z37").Value = -115
$wsFBS.Range("AB38").Value = 2.5
$wsFBS.Range("AF38").Value = -1
$wsFBS.Range("AB39").Value = -9
$wsFBS.Range("AF39").Value = -0.5
$wsFBS.Range("Q40").Value = "ENE"
$wsFBS.Range("Y40").Value = 44.5
$wsFBS.Range("Z40").Value = -110
$wsFBS.Range("AB40").Value = -13.5
$wsFBS.Range("AE40").Value = 0.02298850574712644
$wsFBS.Range("AF40").Value = 0.5
$wsFBS.Range("AB41").Value = -1.5
$wsFBS.Range("AF41").Value = 0.5
$wsFBS.Range("AB44").Value = -26.5
$wsFBS.Range("AF44").Value = 0.5
$wsFBS.Range("AB45").Value = 13.5
$wsFBS.Range("AF45").Value = 0.5
$wsFBS.Range("Y47").Value = 55.5
$wsFBS.Range("AE47").Value = 0.01834862385321101
$wsFBS.Range("Y48").Value = 64.5
$wsFBS.Range("Z48").Value = -105
$wsFBS.Range("AE48").Value = 0.01574803149606299

# --- wind_dir_fg correction (sheet "Other") ---
$wsOther.Range("S38").Value = "ENE"

# --- Timestamp update: every row's Timestamp (col AK, rows 2-51) points to one
#     shared string, so all rows must be rewritten to move that shared text. ---
$newTimestamp = "2024-11-11T10:01:41.103846"
for ($r = 2; $r -le 51; $r++) {
    $wsFBS.Range("AK$r").Value = $newTimestamp
}
